# Actualización automática hashcode - update hash values in column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B15";  Value = "0a25700c1eee3b24b046755a54edb6c1" }
    @{ Cell = "B159"; Value = "e180276079263c04640119ac2f9a2356" }
    @{ Cell = "B169"; Value = "57c8ebb0b1bfe05484cfbeee6e543676" }
    @{ Cell = "B222"; Value = "cd16c0d8eb5dedea53fb77191195f41d" }
    @{ Cell = "B229"; Value = "581c629f37ca65ecb404c64e93a88bc1" }
    @{ Cell = "B339"; Value = "099ff95134ac2a6dda1c6112387b1c53" }
    @{ Cell = "B465"; Value = "14c3403d559313741dc207f2a6e5b139" }
    @{ Cell = "B506"; Value = "d11e0cd41977733cd7b40226af342944" }
    @{ Cell = "B507"; Value = "28b7f4082aa807fa960d3091d6953006" }
    @{ Cell = "B508"; Value = "15a7f9aaaa40c054241246863e869e1e" }
    @{ Cell = "B523"; Value = "9467d99e82d51b4bd4c05a756ee762f7" }
    @{ Cell = "B524"; Value = "b442e64966200cb4be835787721f9bae" }
    @{ Cell = "B555"; Value = "4e3bbd7420307877e649f76fc59c772e" }
    @{ Cell = "B624"; Value = "8eed330081db7ea415c2ac50c2458014" }
    @{ Cell = "B635"; Value = "d450c3da6f90944d2dbd85eeeee6c17e" }
    @{ Cell = "B657"; Value = "107fb1358ef783dcf2c5accf5bcbe860" }
    @{ Cell = "B663"; Value = "576ef50610bc1389eb9c855d927348c8" }
    @{ Cell = "B673"; Value = "5497ef41fdfa9e27c523769b05e1c449" }
    @{ Cell = "B741"; Value = "ecd8d3307557e22a2ca5fe9b25cb8010" }
    @{ Cell = "B827"; Value = "6f14a86add7ba4c658e6672d743c2b75" }
    @{ Cell = "B843"; Value = "08ec81e9257330f99b6ec686fc7b6d56" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
